$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.737058520317078
$ws.Range("B1").Value = 3.80322003364563
$ws.Range("C1").Value = 2.095124244689941
$ws.Range("D1").Value = 1.484809756278992
$ws.Range("E1").Value = 1.268623471260071
